$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new date column header in BB1, copying BA1's formatting (date style)
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)
$ws.Range("BB1").Value2 = 45986

# Carry forward the last known value of each data row into the new BB column
$ws.Range("BB3").Value2 = $ws.Range("BA3").Value2
$ws.Range("BB4").Value2 = $ws.Range("BA4").Value2
$ws.Range("BB5").Value2 = $ws.Range("BA5").Value2
$ws.Range("BB6").Value2 = $ws.Range("BA6").Value2
$ws.Range("BB7").Value2 = $ws.Range("BA7").Value2
$ws.Range("BB8").Value2 = $ws.Range("BA8").Value2
$ws.Range("BB9").Value2 = $ws.Range("BA9").Value2
$ws.Range("BB10").Value2 = $ws.Range("BA10").Value2
$ws.Range("BB11").Value2 = $ws.Range("BA11").Value2
$ws.Range("BB12").Value2 = $ws.Range("BA12").Value2
$ws.Range("BB13").Value2 = $ws.Range("BA13").Value2
$ws.Range("BB14").Value2 = $ws.Range("BA14").Value2
$ws.Range("BB15").Value2 = $ws.Range("BA15").Value2
$ws.Range("BB16").Value2 = $ws.Range("BA16").Value2
$ws.Range("BB17").Value2 = $ws.Range("BA17").Value2
$ws.Range("BB18").Value2 = $ws.Range("BA18").Value2
$ws.Range("BB19").Value2 = $ws.Range("BA19").Value2
$ws.Range("BB20").Value2 = $ws.Range("BA20").Value2
$ws.Range("BB21").Value2 = $ws.Range("BA21").Value2
